# 106年期末績效檢視: update sheet1 ("職能") row 9 with a new self-review
# narrative, clear the now-redundant C9 text, resize the row, and make
# sheet1 the active tab/selection (was previously on sheet2).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "職能"

# Row 9 grows taller to fit the new paragraph.
$ws1.Rows.Item(9).RowHeight = 165

# C9 previously duplicated text that now lives (rewritten) in A9; clear it.
$ws1.Range("C9").Value = ""

# A9 gets the new write-up describing the 預算實支追蹤表單 / 業務員網站
# development work.
$ws1.Range("A9").Value = "    開發預算實支追蹤表單及業務員網站匯款查詢時，因是將資料拋到行政表單系統與業務員網站做呈現，因此需要與數資部的同仁們做溝通，了解他們的程式架構與業務邏輯，在討論的過程中也發生許多問題，像是行政表單系統檔案匯出時裡面包含逗點直接換行導致費用這端在處理資料匯入時會抓到錯誤斷行的資料。與數資部的同仁討論該問題後，提出換行符號變更為特殊符號後，得以解決此問題。因預算實支表單只有每季10號才會使用批次的功能，當批次內容發生錯誤的時候也無法即時重新跑批次產生檔案，只能委託數資部的同仁用資料異動幫忙修改檔案內容，同時在費用系統端也一併對程式做修正，以利下回批次作業可以順利完成。   `n"

# The workbook now opens on sheet1 ("職能") instead of sheet2, with C10
# selected.
$ws1.Activate() | Out-Null
$ws1.Range("C10").Select() | Out-Null
